# Separated controllers into different files
# - Adds row 51 to the AngularJS sheet documenting multi-file controller definition.
# - Adds a new "Webservice" sheet (REST / idempotency notes) before "References".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. AngularJS sheet: append row 51 (Define Controllers in multiple files)
# ---------------------------------------------------------------------------
$angular = $wb.Worksheets.Item("AngularJS")

$angular.Range("A51").Value = "Define Controllers in multiple files"

$controllerNote = "Define module in parent js file as`nvar module = angular.module('app',[]);`nmodule.controller('controller1', function(`$scope){});`nGet the above module in the child js files as:`nvar module = angular.module('app'); // [] needs to be omitted here to get the module declared in above file`nmodule.controller('controller2', function(`$scope){});"
$angular.Range("B51").Value = $controllerNote
$angular.Range("B51").WrapText = $true
$angular.Rows.Item(51).RowHeight = 105

$angular.Activate()
$angular.Range("B57").Select()

# ---------------------------------------------------------------------------
# 2. New "Webservice" sheet, inserted right before "References"
# ---------------------------------------------------------------------------
$referencesSheet = $wb.Worksheets.Item("References")
$ws = $wb.Worksheets.Add($referencesSheet)
$ws.Name = "Webservice"

$ws.Range("A1").Value = "Item"
$ws.Range("B1").Value = "Description"

$ws.Range("A2").Value = "REST"
$ws.Range("B2").Value = "Representational State Transfer"

$ws.Range("A3").Value = "Architecture"
$ws.Range("B3").Value = "REST is an architectural style which is based on web-standards and the HTTP protocol`nEverything is a resource. A resource is accessed via a common interface based on the HTTP standard methods."
$ws.Range("B3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 45

$ws.Range("A4").Value = "REST methods"
$ws.Range("B4").Value = " PUT, GET, POST and DELETE"
$ws.Range("B4").WrapText = $true
$methodsChars = $ws.Range("B4").Characters(2, 26)
$methodsChars.Font.Name = "Times New Roman"
$methodsChars.Font.Size = 13
$ws.Rows.Item(4).RowHeight = 16.5

$ws.Range("A5").Value = "Idempotency"
$ws.Range("B5").Value = "An idempotent method means that the result of a successful performed request is independent of the number of times it is executed. e.g get always reads a resource.Puts always creates a new resource. So both get,  put  and delete requests are idempotent. `nPost is non idempotent`n"
$ws.Range("B5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 75

$ws.Range("A6").Value = "java.lang.IllegalStateException: InjectionManagerFactory not found."
$ws.Range("A6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 30

$ws.Range("A1:B1").WrapText = $true

$ws.Columns.Item(1).ColumnWidth = 60.5703125
$ws.Columns.Item(2).ColumnWidth = 96.5703125

$ws.Activate()
$ws.Range("A10:A11").Select()
